$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("E10").Value = "  +4.26%  "
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +9.50%  "
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("E41").Value = "  +7.95%  "
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("E51").Value = "  -2.44%  "
